# betting_tracker.xlsx update
# - Insert a new "betting_line" column at D on the Bets sheet (goalie_id,
#   team_abbrev, opponent_team, is_home shift right D->E, E->F, F->G, G->H)
# - Repopulate the data rows with the latest scrape: some games still have
#   TBD goalies (cleared goalie/model columns), others got their betting
#   line + a couple of corrected goalie ids / recalculated model outputs.
# - Bump the "Generated:" timestamp on the Summary sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bets")

# --- Header row: shift goalie_id/team_abbrev/opponent_team/is_home right
#     by one and introduce betting_line in column D ---
$ws.Range("D1").Value = "betting_line"
$ws.Range("E1").Value = "goalie_id"
$ws.Range("F1").Value = "team_abbrev"
$ws.Range("G1").Value = "opponent_team"
$ws.Range("H1").Value = "is_home"

# --- Row 2: DAL @ MTL, goalie still TBD ---
$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("F2").Value = "DAL"
$ws.Range("G2").Value = "MTL"
$ws.Range("H2").Value = 1

# --- Row 3: MTL @ DAL, goalie still TBD ---
$ws.Range("C3").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("F3").Value = "MTL"
$ws.Range("G3").Value = "DAL"
$ws.Range("H3").Value = 0

# --- Row 4: CBJ @ PIT, goalie still TBD ---
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = "CBJ"
$ws.Range("G4").Value = "PIT"
$ws.Range("H4").Value = 1

# --- Row 5: PIT @ CBJ, goalie still TBD ---
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("F5").Value = "PIT"
$ws.Range("G5").Value = "CBJ"
$ws.Range("H5").Value = 0

# --- Row 6: FLA @ COL - Tarasov no longer confirmed, back to TBD ---
$ws.Range("C6").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("F6").Value = "FLA"
$ws.Range("G6").Value = "COL"
$ws.Range("H6").Value = 1
$ws.Range("I6").ClearContents()
$ws.Range("J6").ClearContents()
$ws.Range("K6").ClearContents()
$ws.Range("L6").ClearContents()
$ws.Range("M6").ClearContents()

# --- Row 7: COL @ FLA, goalie still TBD ---
$ws.Range("C7").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("F7").Value = "COL"
$ws.Range("G7").Value = "FLA"
$ws.Range("H7").Value = 0
$ws.Range("I7").ClearContents()
$ws.Range("J7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()

# --- Row 8: NJD @ CAR, Allen confirmed - betting_line now known (25.5),
#     predicted_saves refreshed to match the line ---
$ws.Range("D8").Value = 25.5
$ws.Range("E8").Value = 8474596
$ws.Range("F8").Value = "NJD"
$ws.Range("G8").Value = "CAR"
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 25.5

# --- Row 9: CAR @ NJD, goalie id corrected to Bussi's real id (8483548),
#     betting_line now known (23.5), model outputs recalculated ---
$ws.Range("D9").Value = 23.5
$ws.Range("E9").Value = 8483548
$ws.Range("F9").Value = "CAR"
$ws.Range("G9").Value = "NJD"
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 23.5
$ws.Range("J9").Value = 0.4905944764614105
$ws.Range("K9").Value = 1.881104707717896

# --- Row 10: CHI @ VGK, Soderblom confirmed - betting_line now known (25.5),
#     predicted_saves refreshed slightly ---
$ws.Range("D10").Value = 25.5
$ws.Range("E10").Value = 8482821
$ws.Range("F10").Value = "CHI"
$ws.Range("G10").Value = "VGK"
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 25.39999961853027

# --- Row 11: VGK @ CHI, goalie still TBD ---
$ws.Range("C11").ClearContents()
$ws.Range("D11").ClearContents()
$ws.Range("E11").ClearContents()
$ws.Range("F11").Value = "VGK"
$ws.Range("G11").Value = "CHI"
$ws.Range("H11").Value = 0
$ws.Range("I11").ClearContents()
$ws.Range("J11").ClearContents()
$ws.Range("K11").ClearContents()
$ws.Range("L11").ClearContents()
$ws.Range("M11").ClearContents()

# --- Summary sheet: bump the "Generated:" timestamp ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = "2026-01-04 16:16:51"
